# H1AR20 BOM — "PCB is done - Generated all output files"
#
# The part-reference designators U2 (row 23) and U3 (row 24) had their
# component data entered in the wrong rows. This swaps the Manufacturer
# Part Number (B), Description (D), Manufacturer (E), Supplier Part
# Number (F) and Octopart link (H) between the two rows so U2 now shows
# the FT230XQ-R (USB bridge) data and U3 shows the LM3940IMP-3.3 (LDO)
# data. Designator (C) and the two untouched columns (A qty, G price)
# stay put. The Octopart hyperlink that used to live on H23 now has to
# live on H24, and the active selection moves along with the user's last
# click (C24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebuild the sheet's hyperlinks, carrying the H23 link over to H24 ---
# (Re-adding every hyperlink in its original order keeps the relationship
# ids / targets identical to what was already in the file; only the
# anchor cell for the LM3940 datasheet link changes, from H23 to H24.)
$ws.Range("A1:I25").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("H2"), "https://octopart.com/c0603c103k5ractu-kemet-133094?r=sp&s=R_iPBxLnSmGqhkU2rIMFpg")
$ws.Hyperlinks.Add($ws.Range("H3"), "https://octopart.com/cc0805kkx7r7bb105-yageo-8376555?r=sp&s=YKl1wwtkROau_X5nniH2ig")
$ws.Hyperlinks.Add($ws.Range("H10"), "https://octopart.com/vlms1300-gs08-vishay-21709201?r=sp&s=_gcP4_q8T1SC6PJQPTQ9yA")
$ws.Hyperlinks.Add($ws.Range("H18"), "https://octopart.com/crcw060310k0jneb-vishay+dale-46603268")
$ws.Hyperlinks.Add($ws.Range("H19"), "https://octopart.com/crcw060368r0fkea-vishay-39811903")
$ws.Hyperlinks.Add($ws.Range("H20"), "https://octopart.com/crcw060327r0fkea-vishay-39833156")
$ws.Hyperlinks.Add($ws.Range("H7"), "http://octopart.com/tajr105m016rnj-avx-1188552")
$ws.Hyperlinks.Add($ws.Range("H8"), "https://octopart.com/c3225x5r0j107m250ac-tdk-25947844?r=spc")
$ws.Hyperlinks.Add($ws.Range("H9"), "https://octopart.com/c3225x5r0j107m250ac-tdk-25947844?r=spc")
$ws.Hyperlinks.Add($ws.Range("H4"), "https://octopart.com/grm21bc81e475ka12l-murata-10331911?r=sp&s=RY3qZSD8T6mG6TC9CHI5qQ")
$ws.Hyperlinks.Add($ws.Range("H5"), "https://octopart.com/10tpu4r7msi-panasonic-29487748")
$ws.Hyperlinks.Add($ws.Range("H13"), "https://octopart.com/67503-1020-molex-766176?r=sp")
$ws.Hyperlinks.Add($ws.Range("H21"), "https://octopart.com/erj-3geyj390v-panasonic-55422095")
$ws.Hyperlinks.Add($ws.Range("H24"), "https://octopart.com/lm3940imp-3.3%2Fnopb-texas+instruments-24823014?r=sp")
$ws.Hyperlinks.Add($ws.Range("H15"), "https://octopart.com/search?q=RC0603FR-07270RL&start=0")

# --- Swap the U2 (row 23) / U3 (row 24) component data ---
$b23 = $ws.Range("B23").Value2
$d23 = $ws.Range("D23").Value2
$e23 = $ws.Range("E23").Value2
$f23 = $ws.Range("F23").Value2
$h23 = $ws.Range("H23").Value2

$b24 = $ws.Range("B24").Value2
$d24 = $ws.Range("D24").Value2
$e24 = $ws.Range("E24").Value2
$f24 = $ws.Range("F24").Value2
$h24 = $ws.Range("H24").Value2

$ws.Range("B23").Value2 = $b24
$ws.Range("D23").Value2 = $d24
$ws.Range("E23").Value2 = $e24
$ws.Range("F23").Value2 = $f24
$ws.Range("H23").Value2 = $h24

$ws.Range("B24").Value2 = $b23
$ws.Range("D24").Value2 = $d23
$ws.Range("E24").Value2 = $e23
$ws.Range("F24").Value2 = $f23
$ws.Range("H24").Value2 = $h23

# --- Leave the selection where the author last clicked ---
$ws.Range("C24").Select()
